$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("Q:Q").Cut()
$ws.Columns("B:B").Insert()
$ws.Columns("O:O").Cut()
$ws.Columns("E:E").Insert()
$ws.Columns("Q:Q").Cut()
$ws.Columns("F:F").Insert()
$ws.Columns("S:S").Cut()
$ws.Columns("I:I").Insert()
$ws.Columns("S:S").Cut()
$ws.Columns("R:R").Insert()
